$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.38460000000003
$ws.Range("A21").Value = -21.2635
$ws.Range("A23").Value = -21.44560000000002
$ws.Range("A25").Value = -22.41640000000002
